$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header values for columns H and I
$ws.Range("H1").Value = "context"
$ws.Range("I1").Value = "type"

# Data rows 2-8
$contexts = @("Byzantine","Byzantine","Byzantine","Byzantine","Byzantine","Byzantine","Byzantine")
$types    = @("ceiling","floor","PEM (burnt)","PEM","ceiling","ceiling","PEM")

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $contexts[$i]
    $ws.Cells.Item($row, 9).Value = $types[$i]
}

# Apply the (same) font explicitly to force Excel to register a dedicated
# style for the new H:I cells, matching the source workbook's new style.
$ws.Range("H1:I8").Font.Name = "Calibri"

# Column widths to match bestFit (closest achievable values given the
# runtime's column-width quantization)
$ws.Columns.Item(8).ColumnWidth = 8.75
$ws.Columns.Item(9).ColumnWidth = 10.6

# Update selection to match diff
$ws.Range("H1:I8").Select()

$wb.Save()
